# Remove the trailing "blank line / page-break / copyright notice" block
# that used to sit right after the "LOT2039: ... (Requisito fraco)"
# paragraph, while leaving the very last blank + page-break paragraphs
# (the ones that were already at the end of the document) untouched.

$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOT2039" requirement line; the three
# paragraphs that immediately follow it (an empty paragraph, an empty
# paragraph carrying a page break, and the copyright paragraph) are the
# ones that must be deleted.
$marker = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "LOT2039:*Requisito fraco*") {
        $marker = $i
        break
    }
}

if ($marker -eq $null) {
    throw "Could not find the 'LOT2039' requirement paragraph"
}

$firstToRemove = $marker + 1
$lastToRemove = $marker + 3

$start = $d.Paragraphs.Item($firstToRemove).Range.Start
$end = $d.Paragraphs.Item($lastToRemove).Range.End

$r = $d.Range($start, $end)
$r.Delete()
